$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between row 16 and row 17.
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2312"

$ws.Range("F16").Value = 60000
$ws.Range("F17").Value = 14000
